$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(1, 1).Value = "LÍNEA 141 - LP1912 - 26/01/2026"
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:47:32"
$ws.Cells.Item(3, 1).Value = "Total filas: 286"
$ws.Cells.Item(5, 1).Value = "Hora_Scrap"
$ws.Cells.Item(5, 2).Value = "Hora_Llegada"
$ws.Cells.Item(5, 3).Value = "Linea"
$ws.Cells.Item(5, 4).Value = "Minutos"
$ws.Cells.Item(5, 5).Value = "Parada"
$ws.Cells.Item(6, 1).Value = "04:45:48"
$ws.Cells.Item(6, 2).Value = "04:45"
$ws.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = "LP1912"
$ws.Cells.Item(7, 1).Value = "04:21:09"
$ws.Cells.Item(7, 2).Value = "04:46"
$ws.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws.Cells.Item(7, 4).Value = 25
$ws.Cells.Item(7, 5).Value = "LP1912"
$ws.Cells.Item(8, 1).Value = "04:45:48"
$ws.Cells.Item(8, 2).Value = "04:53"
$ws.Cells.Item(8, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(8, 4).Value = 8
$ws.Cells.Item(8, 5).Value = "LP1912"
$ws.Cells.Item(9, 1).Value = "04:56:30"
$ws.Cells.Item(9, 2).Value = "05:16"
$ws.Cells.Item(9, 3).Value = "17_ROMERO"
$ws.Cells.Item(9, 4).Value = 20
$ws.Cells.Item(9, 5).Value = "LP1912"
$ws.Cells.Item(10, 1).Value = "04:56:30"
$ws.Cells.Item(10, 2).Value = "05:22"
$ws.Cells.Item(10, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(10, 4).Value = 26
$ws.Cells.Item(10, 5).Value = "LP1912"
$ws.Cells.Item(11, 1).Value = "05:24:16"
$ws.Cells.Item(11, 2).Value = "05:25"
$ws.Cells.Item(11, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = "LP1912"
$ws.Cells.Item(12, 1).Value = "04:56:30"
$ws.Cells.Item(12, 2).Value = "05:34"
$ws.Cells.Item(12, 3).Value = "215B_EL PATO"
$ws.Cells.Item(12, 4).Value = 38
$ws.Cells.Item(12, 5).Value = "LP1912"
$ws.Cells.Item(13, 1).Value = "05:24:16"
$ws.Cells.Item(13, 2).Value = "05:35"
$ws.Cells.Item(13, 3).Value = "215B_EL PATO"
$ws.Cells.Item(13, 4).Value = 11
$ws.Cells.Item(13, 5).Value = "LP1912"
$ws.Cells.Item(14, 1).Value = "05:24:16"
$ws.Cells.Item(14, 2).Value = "05:46"
$ws.Cells.Item(14, 3).Value = "15_ABASTO"
$ws.Cells.Item(14, 4).Value = 22
$ws.Cells.Item(14, 5).Value = "LP1912"
$ws.Cells.Item(15, 1).Value = "05:24:16"
$ws.Cells.Item(15, 2).Value = "05:54"
$ws.Cells.Item(15, 3).Value = "10_OLMOS"
$ws.Cells.Item(15, 4).Value = 30
$ws.Cells.Item(15, 5).Value = "LP1912"
$ws.Cells.Item(16, 1).Value = "05:55:02"
$ws.Cells.Item(16, 2).Value = "05:55"
$ws.Cells.Item(16, 3).Value = "10_OLMOS"
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = "LP1912"
$ws.Cells.Item(17, 1).Value = "05:24:16"
$ws.Cells.Item(17, 2).Value = "06:04"
$ws.Cells.Item(17, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(17, 4).Value = 40
$ws.Cells.Item(17, 5).Value = "LP1912"
$ws.Cells.Item(18, 1).Value = "05:24:16"
$ws.Cells.Item(18, 2).Value = "06:11"
$ws.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws.Cells.Item(18, 4).Value = 47
$ws.Cells.Item(18, 5).Value = "LP1912"
$ws.Cells.Item(19, 1).Value = "05:55:02"
$ws.Cells.Item(19, 2).Value = "06:12"
$ws.Cells.Item(19, 3).Value = "215A_EL PATO"
$ws.Cells.Item(19, 4).Value = 17
$ws.Cells.Item(19, 5).Value = "LP1912"
$ws.Cells.Item(20, 1).Value = "05:55:02"
$ws.Cells.Item(20, 2).Value = "06:14"
$ws.Cells.Item(20, 3).Value = "225_HARAS DEL SUR"
$ws.Cells.Item(20, 4).Value = 19
$ws.Cells.Item(20, 5).Value = "LP1912"
$ws.Cells.Item(21, 1).Value = "04:56:30"
$ws.Cells.Item(21, 2).Value = "06:18"
$ws.Cells.Item(21, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(21, 4).Value = 82
$ws.Cells.Item(21, 5).Value = "LP1912"
$ws.Cells.Item(22, 1).Value = "05:55:02"
$ws.Cells.Item(22, 2).Value = "06:21"
$ws.Cells.Item(22, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(22, 4).Value = 26
$ws.Cells.Item(22, 5).Value = "LP1912"
$ws.Cells.Item(23, 1).Value = "04:45:48"
$ws.Cells.Item(23, 2).Value = "06:24"
$ws.Cells.Item(23, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(23, 4).Value = 99
$ws.Cells.Item(23, 5).Value = "LP1912"
$ws.Cells.Item(24, 1).Value = "05:55:02"
$ws.Cells.Item(24, 2).Value = "06:27"
$ws.Cells.Item(24, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(24, 4).Value = 32
$ws.Cells.Item(24, 5).Value = "LP1912"
$ws.Cells.Item(25, 1).Value = "06:25:28"
$ws.Cells.Item(25, 2).Value = "06:29"
$ws.Cells.Item(25, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(25, 4).Value = 4
$ws.Cells.Item(25, 5).Value = "LP1912"
$ws.Cells.Item(26, 1).Value = "06:25:28"
$ws.Cells.Item(26, 2).Value = "06:30"
$ws.Cells.Item(26, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(26, 4).Value = 5
$ws.Cells.Item(26, 5).Value = "LP1912"
$ws.Cells.Item(27, 1).Value = "05:55:02"
$ws.Cells.Item(27, 2).Value = "06:30"
$ws.Cells.Item(27, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(27, 4).Value = 35
$ws.Cells.Item(27, 5).Value = "LP1912"
$ws.Cells.Item(28, 1).Value = "05:55:02"
$ws.Cells.Item(28, 2).Value = "06:31"
$ws.Cells.Item(28, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(28, 4).Value = 36
$ws.Cells.Item(28, 5).Value = "LP1912"
$ws.Cells.Item(29, 1).Value = "06:25:28"
$ws.Cells.Item(29, 2).Value = "06:44"
$ws.Cells.Item(29, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(29, 4).Value = 19
$ws.Cells.Item(29, 5).Value = "LP1912"
$ws.Cells.Item(30, 1).Value = "06:25:28"
$ws.Cells.Item(30, 2).Value = "06:46"
$ws.Cells.Item(30, 3).Value = "215C_EL PATO"
$ws.Cells.Item(30, 4).Value = 21
$ws.Cells.Item(30, 5).Value = "LP1912"
$ws.Cells.Item(31, 1).Value = "05:55:02"
$ws.Cells.Item(31, 2).Value = "06:47"
$ws.Cells.Item(31, 3).Value = "215C_EL PATO"
$ws.Cells.Item(31, 4).Value = 52
$ws.Cells.Item(31, 5).Value = "LP1912"
$ws.Cells.Item(32, 1).Value = "06:54:06"
$ws.Cells.Item(32, 2).Value = "06:55"
$ws.Cells.Item(32, 3).Value = "14_ABASTO"
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(32, 5).Value = "LP1912"
$ws.Cells.Item(33, 1).Value = "06:54:06"
$ws.Cells.Item(33, 2).Value = "06:55"
$ws.Cells.Item(33, 3).Value = "215C_EL PATO"
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(33, 5).Value = "LP1912"
$ws.Cells.Item(34, 1).Value = "06:25:28"
$ws.Cells.Item(34, 2).Value = "06:59"
$ws.Cells.Item(34, 3).Value = "14_ABASTO"
$ws.Cells.Item(34, 4).Value = 34
$ws.Cells.Item(34, 5).Value = "LP1912"
$ws.Cells.Item(35, 1).Value = "05:55:02"
$ws.Cells.Item(35, 2).Value = "07:00"
$ws.Cells.Item(35, 3).Value = "14_ABASTO"
$ws.Cells.Item(35, 4).Value = 65
$ws.Cells.Item(35, 5).Value = "LP1912"
$ws.Cells.Item(36, 1).Value = "06:54:06"
$ws.Cells.Item(36, 2).Value = "07:01"
$ws.Cells.Item(36, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(36, 4).Value = 7
$ws.Cells.Item(36, 5).Value = "LP1912"
$ws.Cells.Item(37, 1).Value = "06:25:28"
$ws.Cells.Item(37, 2).Value = "07:05"
$ws.Cells.Item(37, 3).Value = "15_ABASTO"
$ws.Cells.Item(37, 4).Value = 40
$ws.Cells.Item(37, 5).Value = "LP1912"
$ws.Cells.Item(38, 1).Value = "06:54:06"
$ws.Cells.Item(38, 2).Value = "07:05"
$ws.Cells.Item(38, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(38, 4).Value = 11
$ws.Cells.Item(38, 5).Value = "LP1912"
$ws.Cells.Item(39, 1).Value = "06:54:06"
$ws.Cells.Item(39, 2).Value = "07:07"
$ws.Cells.Item(39, 3).Value = "225_GOMEZ"
$ws.Cells.Item(39, 4).Value = 13
$ws.Cells.Item(39, 5).Value = "LP1912"
$ws.Cells.Item(40, 1).Value = "06:54:06"
$ws.Cells.Item(40, 2).Value = "07:07"
$ws.Cells.Item(40, 3).Value = "15_ABASTO"
$ws.Cells.Item(40, 4).Value = 13
$ws.Cells.Item(40, 5).Value = "LP1912"
$ws.Cells.Item(41, 1).Value = "06:25:28"
$ws.Cells.Item(41, 2).Value = "07:11"
$ws.Cells.Item(41, 3).Value = "215A_EL PATO"
$ws.Cells.Item(41, 4).Value = 46
$ws.Cells.Item(41, 5).Value = "LP1912"
$ws.Cells.Item(42, 1).Value = "06:54:06"
$ws.Cells.Item(42, 2).Value = "07:12"
$ws.Cells.Item(42, 3).Value = "215A_EL PATO"
$ws.Cells.Item(42, 4).Value = 18
$ws.Cells.Item(42, 5).Value = "LP1912"
$ws.Cells.Item(43, 1).Value = "06:25:28"
$ws.Cells.Item(43, 2).Value = "07:15"
$ws.Cells.Item(43, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(43, 4).Value = 50
$ws.Cells.Item(43, 5).Value = "LP1912"
$ws.Cells.Item(44, 1).Value = "06:54:06"
$ws.Cells.Item(44, 2).Value = "07:16"
$ws.Cells.Item(44, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(44, 4).Value = 22
$ws.Cells.Item(44, 5).Value = "LP1912"
$ws.Cells.Item(45, 1).Value = "06:54:06"
$ws.Cells.Item(45, 2).Value = "07:17"
$ws.Cells.Item(45, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(45, 4).Value = 23
$ws.Cells.Item(45, 5).Value = "LP1912"
$ws.Cells.Item(46, 1).Value = "07:17:59"
$ws.Cells.Item(46, 2).Value = "07:20"
$ws.Cells.Item(46, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(46, 4).Value = 3
$ws.Cells.Item(46, 5).Value = "LP1912"
$ws.Cells.Item(47, 1).Value = "06:54:06"
$ws.Cells.Item(47, 2).Value = "07:21"
$ws.Cells.Item(47, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(47, 4).Value = 27
$ws.Cells.Item(47, 5).Value = "LP1912"
$ws.Cells.Item(48, 1).Value = "06:54:06"
$ws.Cells.Item(48, 2).Value = "07:23"
$ws.Cells.Item(48, 3).Value = "10_OLMOS"
$ws.Cells.Item(48, 4).Value = 29
$ws.Cells.Item(48, 5).Value = "LP1912"
$ws.Cells.Item(49, 1).Value = "07:17:59"
$ws.Cells.Item(49, 2).Value = "07:31"
$ws.Cells.Item(49, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(49, 4).Value = 14
$ws.Cells.Item(49, 5).Value = "LP1912"
$ws.Cells.Item(50, 1).Value = "07:17:59"
$ws.Cells.Item(50, 2).Value = "07:31"
$ws.Cells.Item(50, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(50, 4).Value = 14
$ws.Cells.Item(50, 5).Value = "LP1912"
$ws.Cells.Item(51, 1).Value = "07:17:59"
$ws.Cells.Item(51, 2).Value = "07:31"
$ws.Cells.Item(51, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(51, 4).Value = 14
$ws.Cells.Item(51, 5).Value = "LP1912"
$ws.Cells.Item(52, 1).Value = "05:55:02"
$ws.Cells.Item(52, 2).Value = "07:32"
$ws.Cells.Item(52, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(52, 4).Value = 97
$ws.Cells.Item(52, 5).Value = "LP1912"
$ws.Cells.Item(53, 1).Value = "06:54:06"
$ws.Cells.Item(53, 2).Value = "07:32"
$ws.Cells.Item(53, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(53, 4).Value = 38
$ws.Cells.Item(53, 5).Value = "LP1912"
$ws.Cells.Item(54, 1).Value = "06:54:06"
$ws.Cells.Item(54, 2).Value = "07:32"
$ws.Cells.Item(54, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(54, 4).Value = 38
$ws.Cells.Item(54, 5).Value = "LP1912"
$ws.Cells.Item(55, 1).Value = "07:17:59"
$ws.Cells.Item(55, 2).Value = "07:35"
$ws.Cells.Item(55, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(55, 4).Value = 18
$ws.Cells.Item(55, 5).Value = "LP1912"
$ws.Cells.Item(56, 1).Value = "07:17:59"
$ws.Cells.Item(56, 2).Value = "07:36"
$ws.Cells.Item(56, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(56, 4).Value = 19
$ws.Cells.Item(56, 5).Value = "LP1912"
$ws.Cells.Item(57, 1).Value = "06:54:06"
$ws.Cells.Item(57, 2).Value = "07:37"
$ws.Cells.Item(57, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(57, 4).Value = 43
$ws.Cells.Item(57, 5).Value = "LP1912"
$ws.Cells.Item(58, 1).Value = "07:17:59"
$ws.Cells.Item(58, 2).Value = "07:38"
$ws.Cells.Item(58, 3).Value = "10_OLMOS"
$ws.Cells.Item(58, 4).Value = 21
$ws.Cells.Item(58, 5).Value = "LP1912"
$ws.Cells.Item(59, 1).Value = "06:54:06"
$ws.Cells.Item(59, 2).Value = "07:39"
$ws.Cells.Item(59, 3).Value = "10_OLMOS"
$ws.Cells.Item(59, 4).Value = 45
$ws.Cells.Item(59, 5).Value = "LP1912"
$ws.Cells.Item(60, 1).Value = "07:17:59"
$ws.Cells.Item(60, 2).Value = "07:46"
$ws.Cells.Item(60, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(60, 4).Value = 29
$ws.Cells.Item(60, 5).Value = "LP1912"
$ws.Cells.Item(61, 1).Value = "07:17:59"
$ws.Cells.Item(61, 2).Value = "07:47"
$ws.Cells.Item(61, 3).Value = "14_ABASTO"
$ws.Cells.Item(61, 4).Value = 30
$ws.Cells.Item(61, 5).Value = "LP1912"
$ws.Cells.Item(62, 1).Value = "07:48:05"
$ws.Cells.Item(62, 2).Value = "07:48"
$ws.Cells.Item(62, 3).Value = "14_ABASTO"
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = "LP1912"
$ws.Cells.Item(63, 1).Value = "07:17:59"
$ws.Cells.Item(63, 2).Value = "07:51"
$ws.Cells.Item(63, 3).Value = "215D_EL PATO"
$ws.Cells.Item(63, 4).Value = 34
$ws.Cells.Item(63, 5).Value = "LP1912"
$ws.Cells.Item(64, 1).Value = "07:48:05"
$ws.Cells.Item(64, 2).Value = "07:52"
$ws.Cells.Item(64, 3).Value = "215D_EL PATO"
$ws.Cells.Item(64, 4).Value = 4
$ws.Cells.Item(64, 5).Value = "LP1912"
$ws.Cells.Item(65, 1).Value = "07:17:59"
$ws.Cells.Item(65, 2).Value = "07:59"
$ws.Cells.Item(65, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(65, 4).Value = 42
$ws.Cells.Item(65, 5).Value = "LP1912"
$ws.Cells.Item(66, 1).Value = "07:48:05"
$ws.Cells.Item(66, 2).Value = "08:02"
$ws.Cells.Item(66, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(66, 4).Value = 14
$ws.Cells.Item(66, 5).Value = "LP1912"
$ws.Cells.Item(67, 1).Value = "08:01:10"
$ws.Cells.Item(67, 2).Value = "08:03"
$ws.Cells.Item(67, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(67, 4).Value = 2
$ws.Cells.Item(67, 5).Value = "LP1912"
$ws.Cells.Item(68, 1).Value = "07:17:59"
$ws.Cells.Item(68, 2).Value = "08:03"
$ws.Cells.Item(68, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(68, 4).Value = 46
$ws.Cells.Item(68, 5).Value = "LP1912"
$ws.Cells.Item(69, 1).Value = "08:01:10"
$ws.Cells.Item(69, 2).Value = "08:04"
$ws.Cells.Item(69, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(69, 4).Value = 3
$ws.Cells.Item(69, 5).Value = "LP1912"
$ws.Cells.Item(70, 1).Value = "06:54:06"
$ws.Cells.Item(70, 2).Value = "08:06"
$ws.Cells.Item(70, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(70, 4).Value = 72
$ws.Cells.Item(70, 5).Value = "LP1912"
$ws.Cells.Item(71, 1).Value = "07:17:59"
$ws.Cells.Item(71, 2).Value = "08:11"
$ws.Cells.Item(71, 3).Value = "15_ABASTO"
$ws.Cells.Item(71, 4).Value = 54
$ws.Cells.Item(71, 5).Value = "LP1912"
$ws.Cells.Item(72, 1).Value = "08:01:10"
$ws.Cells.Item(72, 2).Value = "08:11"
$ws.Cells.Item(72, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(72, 4).Value = 10
$ws.Cells.Item(72, 5).Value = "LP1912"
$ws.Cells.Item(73, 1).Value = "08:01:10"
$ws.Cells.Item(73, 2).Value = "08:12"
$ws.Cells.Item(73, 3).Value = "15_ABASTO"
$ws.Cells.Item(73, 4).Value = 11
$ws.Cells.Item(73, 5).Value = "LP1912"
$ws.Cells.Item(74, 1).Value = "07:17:59"
$ws.Cells.Item(74, 2).Value = "08:12"
$ws.Cells.Item(74, 3).Value = "10_OLMOS"
$ws.Cells.Item(74, 4).Value = 55
$ws.Cells.Item(74, 5).Value = "LP1912"
$ws.Cells.Item(75, 1).Value = "08:01:10"
$ws.Cells.Item(75, 2).Value = "08:13"
$ws.Cells.Item(75, 3).Value = "10_OLMOS"
$ws.Cells.Item(75, 4).Value = 12
$ws.Cells.Item(75, 5).Value = "LP1912"
$ws.Cells.Item(76, 1).Value = "07:17:59"
$ws.Cells.Item(76, 2).Value = "08:20"
$ws.Cells.Item(76, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(76, 4).Value = 63
$ws.Cells.Item(76, 5).Value = "LP1912"
$ws.Cells.Item(77, 1).Value = "06:54:06"
$ws.Cells.Item(77, 2).Value = "08:21"
$ws.Cells.Item(77, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(77, 4).Value = 87
$ws.Cells.Item(77, 5).Value = "LP1912"
$ws.Cells.Item(78, 1).Value = "07:17:59"
$ws.Cells.Item(78, 2).Value = "08:22"
$ws.Cells.Item(78, 3).Value = "215B_EL PATO"
$ws.Cells.Item(78, 4).Value = 65
$ws.Cells.Item(78, 5).Value = "LP1912"
$ws.Cells.Item(79, 1).Value = "07:17:59"
$ws.Cells.Item(79, 2).Value = "08:22"
$ws.Cells.Item(79, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(79, 4).Value = 65
$ws.Cells.Item(79, 5).Value = "LP1912"
$ws.Cells.Item(80, 1).Value = "08:01:10"
$ws.Cells.Item(80, 2).Value = "08:23"
$ws.Cells.Item(80, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(80, 4).Value = 22
$ws.Cells.Item(80, 5).Value = "LP1912"
$ws.Cells.Item(81, 1).Value = "08:01:10"
$ws.Cells.Item(81, 2).Value = "08:23"
$ws.Cells.Item(81, 3).Value = "215B_EL PATO"
$ws.Cells.Item(81, 4).Value = 22
$ws.Cells.Item(81, 5).Value = "LP1912"
$ws.Cells.Item(82, 1).Value = "08:01:10"
$ws.Cells.Item(82, 2).Value = "08:24"
$ws.Cells.Item(82, 3).Value = "14_ABASTO"
$ws.Cells.Item(82, 4).Value = 23
$ws.Cells.Item(82, 5).Value = "LP1912"
$ws.Cells.Item(83, 1).Value = "07:17:59"
$ws.Cells.Item(83, 2).Value = "08:26"
$ws.Cells.Item(83, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(83, 4).Value = 69
$ws.Cells.Item(83, 5).Value = "LP1912"
$ws.Cells.Item(84, 1).Value = "08:01:10"
$ws.Cells.Item(84, 2).Value = "08:27"
$ws.Cells.Item(84, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(84, 4).Value = 26
$ws.Cells.Item(84, 5).Value = "LP1912"
$ws.Cells.Item(85, 1).Value = "07:48:05"
$ws.Cells.Item(85, 2).Value = "08:30"
$ws.Cells.Item(85, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(85, 4).Value = 42
$ws.Cells.Item(85, 5).Value = "LP1912"
$ws.Cells.Item(86, 1).Value = "08:01:10"
$ws.Cells.Item(86, 2).Value = "08:33"
$ws.Cells.Item(86, 3).Value = "10_OLMOS"
$ws.Cells.Item(86, 4).Value = 32
$ws.Cells.Item(86, 5).Value = "LP1912"
$ws.Cells.Item(87, 1).Value = "08:31:01"
$ws.Cells.Item(87, 2).Value = "08:35"
$ws.Cells.Item(87, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(87, 4).Value = 4
$ws.Cells.Item(87, 5).Value = "LP1912"
$ws.Cells.Item(88, 1).Value = "07:48:05"
$ws.Cells.Item(88, 2).Value = "08:37"
$ws.Cells.Item(88, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(88, 4).Value = 49
$ws.Cells.Item(88, 5).Value = "LP1912"
$ws.Cells.Item(89, 1).Value = "08:01:10"
$ws.Cells.Item(89, 2).Value = "08:40"
$ws.Cells.Item(89, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(89, 4).Value = 39
$ws.Cells.Item(89, 5).Value = "LP1912"
$ws.Cells.Item(90, 1).Value = "07:17:59"
$ws.Cells.Item(90, 2).Value = "08:41"
$ws.Cells.Item(90, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(90, 4).Value = 84
$ws.Cells.Item(90, 5).Value = "LP1912"
$ws.Cells.Item(91, 1).Value = "08:31:01"
$ws.Cells.Item(91, 2).Value = "08:42"
$ws.Cells.Item(91, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(91, 4).Value = 11
$ws.Cells.Item(91, 5).Value = "LP1912"
$ws.Cells.Item(92, 1).Value = "07:17:59"
$ws.Cells.Item(92, 2).Value = "08:43"
$ws.Cells.Item(92, 3).Value = "14_ABASTO"
$ws.Cells.Item(92, 4).Value = 86
$ws.Cells.Item(92, 5).Value = "LP1912"
$ws.Cells.Item(93, 1).Value = "08:31:01"
$ws.Cells.Item(93, 2).Value = "08:44"
$ws.Cells.Item(93, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(93, 4).Value = 13
$ws.Cells.Item(93, 5).Value = "LP1912"
$ws.Cells.Item(94, 1).Value = "08:31:01"
$ws.Cells.Item(94, 2).Value = "08:44"
$ws.Cells.Item(94, 3).Value = "14_ABASTO"
$ws.Cells.Item(94, 4).Value = 13
$ws.Cells.Item(94, 5).Value = "LP1912"
$ws.Cells.Item(95, 1).Value = "08:47:26"
$ws.Cells.Item(95, 2).Value = "08:47"
$ws.Cells.Item(95, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(95, 4).Value = 0
$ws.Cells.Item(95, 5).Value = "LP1912"
$ws.Cells.Item(96, 1).Value = "08:01:10"
$ws.Cells.Item(96, 2).Value = "08:49"
$ws.Cells.Item(96, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(96, 4).Value = 48
$ws.Cells.Item(96, 5).Value = "LP1912"
$ws.Cells.Item(97, 1).Value = "08:47:26"
$ws.Cells.Item(97, 2).Value = "08:53"
$ws.Cells.Item(97, 3).Value = "10_OLMOS"
$ws.Cells.Item(97, 4).Value = 6
$ws.Cells.Item(97, 5).Value = "LP1912"
$ws.Cells.Item(98, 1).Value = "07:17:59"
$ws.Cells.Item(98, 2).Value = "08:53"
$ws.Cells.Item(98, 3).Value = "17_ROMERO"
$ws.Cells.Item(98, 4).Value = 96
$ws.Cells.Item(98, 5).Value = "LP1912"
$ws.Cells.Item(99, 1).Value = "08:47:26"
$ws.Cells.Item(99, 2).Value = "08:54"
$ws.Cells.Item(99, 3).Value = "17_ROMERO"
$ws.Cells.Item(99, 4).Value = 7
$ws.Cells.Item(99, 5).Value = "LP1912"
$ws.Cells.Item(100, 1).Value = "08:55:01"
$ws.Cells.Item(100, 2).Value = "08:56"
$ws.Cells.Item(100, 3).Value = "17_ROMERO"
$ws.Cells.Item(100, 4).Value = 1
$ws.Cells.Item(100, 5).Value = "LP1912"
$ws.Cells.Item(101, 1).Value = "08:47:26"
$ws.Cells.Item(101, 2).Value = "09:01"
$ws.Cells.Item(101, 3).Value = "215A_EL PATO"
$ws.Cells.Item(101, 4).Value = 14
$ws.Cells.Item(101, 5).Value = "LP1912"
$ws.Cells.Item(102, 1).Value = "08:55:01"
$ws.Cells.Item(102, 2).Value = "09:02"
$ws.Cells.Item(102, 3).Value = "215A_EL PATO"
$ws.Cells.Item(102, 4).Value = 7
$ws.Cells.Item(102, 5).Value = "LP1912"
$ws.Cells.Item(103, 1).Value = "08:47:26"
$ws.Cells.Item(103, 2).Value = "09:03"
$ws.Cells.Item(103, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(103, 4).Value = 16
$ws.Cells.Item(103, 5).Value = "LP1912"
$ws.Cells.Item(104, 1).Value = "08:55:01"
$ws.Cells.Item(104, 2).Value = "09:04"
$ws.Cells.Item(104, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(104, 4).Value = 9
$ws.Cells.Item(104, 5).Value = "LP1912"
$ws.Cells.Item(105, 1).Value = "08:47:26"
$ws.Cells.Item(105, 2).Value = "09:05"
$ws.Cells.Item(105, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(105, 4).Value = 18
$ws.Cells.Item(105, 5).Value = "LP1912"
$ws.Cells.Item(106, 1).Value = "08:55:01"
$ws.Cells.Item(106, 2).Value = "09:06"
$ws.Cells.Item(106, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(106, 4).Value = 11
$ws.Cells.Item(106, 5).Value = "LP1912"
$ws.Cells.Item(107, 1).Value = "08:47:26"
$ws.Cells.Item(107, 2).Value = "09:10"
$ws.Cells.Item(107, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(107, 4).Value = 23
$ws.Cells.Item(107, 5).Value = "LP1912"
$ws.Cells.Item(108, 1).Value = "08:55:01"
$ws.Cells.Item(108, 2).Value = "09:11"
$ws.Cells.Item(108, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(108, 4).Value = 16
$ws.Cells.Item(108, 5).Value = "LP1912"
$ws.Cells.Item(109, 1).Value = "08:47:26"
$ws.Cells.Item(109, 2).Value = "09:16"
$ws.Cells.Item(109, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(109, 4).Value = 29
$ws.Cells.Item(109, 5).Value = "LP1912"
$ws.Cells.Item(110, 1).Value = "08:55:01"
$ws.Cells.Item(110, 2).Value = "09:17"
$ws.Cells.Item(110, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(110, 4).Value = 22
$ws.Cells.Item(110, 5).Value = "LP1912"
$ws.Cells.Item(111, 1).Value = "08:31:01"
$ws.Cells.Item(111, 2).Value = "09:19"
$ws.Cells.Item(111, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(111, 4).Value = 48
$ws.Cells.Item(111, 5).Value = "LP1912"
$ws.Cells.Item(112, 1).Value = "08:55:01"
$ws.Cells.Item(112, 2).Value = "09:21"
$ws.Cells.Item(112, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(112, 4).Value = 26
$ws.Cells.Item(112, 5).Value = "LP1912"
$ws.Cells.Item(113, 1).Value = "08:47:26"
$ws.Cells.Item(113, 2).Value = "09:22"
$ws.Cells.Item(113, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(113, 4).Value = 35
$ws.Cells.Item(113, 5).Value = "LP1912"
$ws.Cells.Item(114, 1).Value = "08:55:01"
$ws.Cells.Item(114, 2).Value = "09:23"
$ws.Cells.Item(114, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(114, 4).Value = 28
$ws.Cells.Item(114, 5).Value = "LP1912"
$ws.Cells.Item(115, 1).Value = "08:47:26"
$ws.Cells.Item(115, 2).Value = "09:23"
$ws.Cells.Item(115, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(115, 4).Value = 36
$ws.Cells.Item(115, 5).Value = "LP1912"
$ws.Cells.Item(116, 1).Value = "08:31:01"
$ws.Cells.Item(116, 2).Value = "09:23"
$ws.Cells.Item(116, 3).Value = "17_ROMERO"
$ws.Cells.Item(116, 4).Value = 52
$ws.Cells.Item(116, 5).Value = "LP1912"
$ws.Cells.Item(117, 1).Value = "08:55:01"
$ws.Cells.Item(117, 2).Value = "09:24"
$ws.Cells.Item(117, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(117, 4).Value = 29
$ws.Cells.Item(117, 5).Value = "LP1912"
$ws.Cells.Item(118, 1).Value = "09:28:08"
$ws.Cells.Item(118, 2).Value = "09:29"
$ws.Cells.Item(118, 3).Value = "17_ROMERO"
$ws.Cells.Item(118, 4).Value = 1
$ws.Cells.Item(118, 5).Value = "LP1912"
$ws.Cells.Item(119, 1).Value = "07:48:05"
$ws.Cells.Item(119, 2).Value = "09:32"
$ws.Cells.Item(119, 3).Value = "15_ABASTO"
$ws.Cells.Item(119, 4).Value = 104
$ws.Cells.Item(119, 5).Value = "LP1912"
$ws.Cells.Item(120, 1).Value = "09:28:08"
$ws.Cells.Item(120, 2).Value = "09:33"
$ws.Cells.Item(120, 3).Value = "10_OLMOS"
$ws.Cells.Item(120, 4).Value = 5
$ws.Cells.Item(120, 5).Value = "LP1912"
$ws.Cells.Item(121, 1).Value = "08:47:26"
$ws.Cells.Item(121, 2).Value = "09:34"
$ws.Cells.Item(121, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(121, 4).Value = 47
$ws.Cells.Item(121, 5).Value = "LP1912"
$ws.Cells.Item(122, 1).Value = "07:48:05"
$ws.Cells.Item(122, 2).Value = "09:34"
$ws.Cells.Item(122, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(122, 4).Value = 106
$ws.Cells.Item(122, 5).Value = "LP1912"
$ws.Cells.Item(123, 1).Value = "09:28:08"
$ws.Cells.Item(123, 2).Value = "09:35"
$ws.Cells.Item(123, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(123, 4).Value = 7
$ws.Cells.Item(123, 5).Value = "LP1912"
$ws.Cells.Item(124, 1).Value = "08:55:01"
$ws.Cells.Item(124, 2).Value = "09:35"
$ws.Cells.Item(124, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(124, 4).Value = 40
$ws.Cells.Item(124, 5).Value = "LP1912"
$ws.Cells.Item(125, 1).Value = "08:47:26"
$ws.Cells.Item(125, 2).Value = "09:35"
$ws.Cells.Item(125, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(125, 4).Value = 48
$ws.Cells.Item(125, 5).Value = "LP1912"
$ws.Cells.Item(126, 1).Value = "09:28:08"
$ws.Cells.Item(126, 2).Value = "09:42"
$ws.Cells.Item(126, 3).Value = "215C_EL PATO"
$ws.Cells.Item(126, 4).Value = 14
$ws.Cells.Item(126, 5).Value = "LP1912"
$ws.Cells.Item(127, 1).Value = "08:55:01"
$ws.Cells.Item(127, 2).Value = "09:43"
$ws.Cells.Item(127, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(127, 4).Value = 48
$ws.Cells.Item(127, 5).Value = "LP1912"
$ws.Cells.Item(128, 1).Value = "08:47:26"
$ws.Cells.Item(128, 2).Value = "09:43"
$ws.Cells.Item(128, 3).Value = "14_ABASTO"
$ws.Cells.Item(128, 4).Value = 56
$ws.Cells.Item(128, 5).Value = "LP1912"
$ws.Cells.Item(129, 1).Value = "09:28:08"
$ws.Cells.Item(129, 2).Value = "09:44"
$ws.Cells.Item(129, 3).Value = "14_ABASTO"
$ws.Cells.Item(129, 4).Value = 16
$ws.Cells.Item(129, 5).Value = "LP1912"
$ws.Cells.Item(130, 1).Value = "09:28:08"
$ws.Cells.Item(130, 2).Value = "09:51"
$ws.Cells.Item(130, 3).Value = "10_OLMOS"
$ws.Cells.Item(130, 4).Value = 23
$ws.Cells.Item(130, 5).Value = "LP1912"
$ws.Cells.Item(131, 1).Value = "09:28:08"
$ws.Cells.Item(131, 2).Value = "09:52"
$ws.Cells.Item(131, 3).Value = "15_ABASTO"
$ws.Cells.Item(131, 4).Value = 24
$ws.Cells.Item(131, 5).Value = "LP1912"
$ws.Cells.Item(132, 1).Value = "09:28:08"
$ws.Cells.Item(132, 2).Value = "09:53"
$ws.Cells.Item(132, 3).Value = "10_OLMOS"
$ws.Cells.Item(132, 4).Value = 25
$ws.Cells.Item(132, 5).Value = "LP1912"
$ws.Cells.Item(133, 1).Value = "08:31:01"
$ws.Cells.Item(133, 2).Value = "09:56"
$ws.Cells.Item(133, 3).Value = "10_OLMOS"
$ws.Cells.Item(133, 4).Value = 85
$ws.Cells.Item(133, 5).Value = "LP1912"
$ws.Cells.Item(134, 1).Value = "09:28:08"
$ws.Cells.Item(134, 2).Value = "10:04"
$ws.Cells.Item(134, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(134, 4).Value = 36
$ws.Cells.Item(134, 5).Value = "LP1912"
$ws.Cells.Item(135, 1).Value = "09:28:08"
$ws.Cells.Item(135, 2).Value = "10:05"
$ws.Cells.Item(135, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(135, 4).Value = 37
$ws.Cells.Item(135, 5).Value = "LP1912"
$ws.Cells.Item(136, 1).Value = "08:47:26"
$ws.Cells.Item(136, 2).Value = "10:10"
$ws.Cells.Item(136, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(136, 4).Value = 83
$ws.Cells.Item(136, 5).Value = "LP1912"
$ws.Cells.Item(137, 1).Value = "09:28:08"
$ws.Cells.Item(137, 2).Value = "10:11"
$ws.Cells.Item(137, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(137, 4).Value = 43
$ws.Cells.Item(137, 5).Value = "LP1912"
$ws.Cells.Item(138, 1).Value = "09:28:08"
$ws.Cells.Item(138, 2).Value = "10:12"
$ws.Cells.Item(138, 3).Value = "15_ABASTO"
$ws.Cells.Item(138, 4).Value = 44
$ws.Cells.Item(138, 5).Value = "LP1912"
$ws.Cells.Item(139, 1).Value = "09:28:08"
$ws.Cells.Item(139, 2).Value = "10:16"
$ws.Cells.Item(139, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(139, 4).Value = 48
$ws.Cells.Item(139, 5).Value = "LP1912"
$ws.Cells.Item(140, 1).Value = "09:28:08"
$ws.Cells.Item(140, 2).Value = "10:21"
$ws.Cells.Item(140, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(140, 4).Value = 53
$ws.Cells.Item(140, 5).Value = "LP1912"
$ws.Cells.Item(141, 1).Value = "09:28:08"
$ws.Cells.Item(141, 2).Value = "10:23"
$ws.Cells.Item(141, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(141, 4).Value = 55
$ws.Cells.Item(141, 5).Value = "LP1912"
$ws.Cells.Item(142, 1).Value = "09:28:08"
$ws.Cells.Item(142, 2).Value = "10:24"
$ws.Cells.Item(142, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(142, 4).Value = 56
$ws.Cells.Item(142, 5).Value = "LP1912"
$ws.Cells.Item(143, 1).Value = "10:25:24"
$ws.Cells.Item(143, 2).Value = "10:26"
$ws.Cells.Item(143, 3).Value = "215A_EL PATO"
$ws.Cells.Item(143, 4).Value = 1
$ws.Cells.Item(143, 5).Value = "LP1912"
$ws.Cells.Item(144, 1).Value = "08:55:01"
$ws.Cells.Item(144, 2).Value = "10:26"
$ws.Cells.Item(144, 3).Value = "10_OLMOS"
$ws.Cells.Item(144, 4).Value = 91
$ws.Cells.Item(144, 5).Value = "LP1912"
$ws.Cells.Item(145, 1).Value = "10:25:24"
$ws.Cells.Item(145, 2).Value = "10:27"
$ws.Cells.Item(145, 3).Value = "17_ROMERO"
$ws.Cells.Item(145, 4).Value = 2
$ws.Cells.Item(145, 5).Value = "LP1912"
$ws.Cells.Item(146, 1).Value = "09:28:08"
$ws.Cells.Item(146, 2).Value = "10:27"
$ws.Cells.Item(146, 3).Value = "215A_EL PATO"
$ws.Cells.Item(146, 4).Value = 59
$ws.Cells.Item(146, 5).Value = "LP1912"
$ws.Cells.Item(147, 1).Value = "10:25:24"
$ws.Cells.Item(147, 2).Value = "10:33"
$ws.Cells.Item(147, 3).Value = "10_OLMOS"
$ws.Cells.Item(147, 4).Value = 8
$ws.Cells.Item(147, 5).Value = "LP1912"
$ws.Cells.Item(148, 1).Value = "10:25:24"
$ws.Cells.Item(148, 2).Value = "10:34"
$ws.Cells.Item(148, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(148, 4).Value = 9
$ws.Cells.Item(148, 5).Value = "LP1912"
$ws.Cells.Item(149, 1).Value = "10:25:24"
$ws.Cells.Item(149, 2).Value = "10:35"
$ws.Cells.Item(149, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(149, 4).Value = 10
$ws.Cells.Item(149, 5).Value = "LP1912"
$ws.Cells.Item(150, 1).Value = "09:28:08"
$ws.Cells.Item(150, 2).Value = "10:35"
$ws.Cells.Item(150, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(150, 4).Value = 67
$ws.Cells.Item(150, 5).Value = "LP1912"
$ws.Cells.Item(151, 1).Value = "10:25:24"
$ws.Cells.Item(151, 2).Value = "10:42"
$ws.Cells.Item(151, 3).Value = "17_ROMERO"
$ws.Cells.Item(151, 4).Value = 17
$ws.Cells.Item(151, 5).Value = "LP1912"
$ws.Cells.Item(152, 1).Value = "10:25:24"
$ws.Cells.Item(152, 2).Value = "10:43"
$ws.Cells.Item(152, 3).Value = "14_ABASTO"
$ws.Cells.Item(152, 4).Value = 18
$ws.Cells.Item(152, 5).Value = "LP1912"
$ws.Cells.Item(153, 1).Value = "09:28:08"
$ws.Cells.Item(153, 2).Value = "10:44"
$ws.Cells.Item(153, 3).Value = "14_ABASTO"
$ws.Cells.Item(153, 4).Value = 76
$ws.Cells.Item(153, 5).Value = "LP1912"
$ws.Cells.Item(154, 1).Value = "10:25:24"
$ws.Cells.Item(154, 2).Value = "10:46"
$ws.Cells.Item(154, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(154, 4).Value = 21
$ws.Cells.Item(154, 5).Value = "LP1912"
$ws.Cells.Item(155, 1).Value = "10:25:24"
$ws.Cells.Item(155, 2).Value = "10:52"
$ws.Cells.Item(155, 3).Value = "15_ABASTO"
$ws.Cells.Item(155, 4).Value = 27
$ws.Cells.Item(155, 5).Value = "LP1912"
$ws.Cells.Item(156, 1).Value = "10:25:24"
$ws.Cells.Item(156, 2).Value = "10:53"
$ws.Cells.Item(156, 3).Value = "10_OLMOS"
$ws.Cells.Item(156, 4).Value = 28
$ws.Cells.Item(156, 5).Value = "LP1912"
$ws.Cells.Item(157, 1).Value = "10:25:24"
$ws.Cells.Item(157, 2).Value = "10:57"
$ws.Cells.Item(157, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(157, 4).Value = 32
$ws.Cells.Item(157, 5).Value = "LP1912"
$ws.Cells.Item(158, 1).Value = "10:58:14"
$ws.Cells.Item(158, 2).Value = "10:58"
$ws.Cells.Item(158, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = "LP1912"
$ws.Cells.Item(159, 1).Value = "09:28:08"
$ws.Cells.Item(159, 2).Value = "11:01"
$ws.Cells.Item(159, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(159, 4).Value = 93
$ws.Cells.Item(159, 5).Value = "LP1912"
$ws.Cells.Item(160, 1).Value = "10:58:14"
$ws.Cells.Item(160, 2).Value = "11:02"
$ws.Cells.Item(160, 3).Value = "215C_EL PATO"
$ws.Cells.Item(160, 4).Value = 4
$ws.Cells.Item(160, 5).Value = "LP1912"
$ws.Cells.Item(161, 1).Value = "10:25:24"
$ws.Cells.Item(161, 2).Value = "11:03"
$ws.Cells.Item(161, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(161, 4).Value = 38
$ws.Cells.Item(161, 5).Value = "LP1912"
$ws.Cells.Item(162, 1).Value = "10:25:24"
$ws.Cells.Item(162, 2).Value = "11:04"
$ws.Cells.Item(162, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(162, 4).Value = 39
$ws.Cells.Item(162, 5).Value = "LP1912"
$ws.Cells.Item(163, 1).Value = "10:58:14"
$ws.Cells.Item(163, 2).Value = "11:04"
$ws.Cells.Item(163, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(163, 4).Value = 6
$ws.Cells.Item(163, 5).Value = "LP1912"
$ws.Cells.Item(164, 1).Value = "10:25:24"
$ws.Cells.Item(164, 2).Value = "11:06"
$ws.Cells.Item(164, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(164, 4).Value = 41
$ws.Cells.Item(164, 5).Value = "LP1912"
$ws.Cells.Item(165, 1).Value = "10:58:14"
$ws.Cells.Item(165, 2).Value = "11:07"
$ws.Cells.Item(165, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(165, 4).Value = 9
$ws.Cells.Item(165, 5).Value = "LP1912"
$ws.Cells.Item(166, 1).Value = "10:58:14"
$ws.Cells.Item(166, 2).Value = "11:07"
$ws.Cells.Item(166, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(166, 4).Value = 9
$ws.Cells.Item(166, 5).Value = "LP1912"
$ws.Cells.Item(167, 1).Value = "10:58:14"
$ws.Cells.Item(167, 2).Value = "11:11"
$ws.Cells.Item(167, 3).Value = "10_OLMOS"
$ws.Cells.Item(167, 4).Value = 13
$ws.Cells.Item(167, 5).Value = "LP1912"
$ws.Cells.Item(168, 1).Value = "10:58:14"
$ws.Cells.Item(168, 2).Value = "11:12"
$ws.Cells.Item(168, 3).Value = "15_ABASTO"
$ws.Cells.Item(168, 4).Value = 14
$ws.Cells.Item(168, 5).Value = "LP1912"
$ws.Cells.Item(169, 1).Value = "10:25:24"
$ws.Cells.Item(169, 2).Value = "11:19"
$ws.Cells.Item(169, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(169, 4).Value = 54
$ws.Cells.Item(169, 5).Value = "LP1912"
$ws.Cells.Item(170, 1).Value = "10:58:14"
$ws.Cells.Item(170, 2).Value = "11:20"
$ws.Cells.Item(170, 3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(170, 4).Value = 22
$ws.Cells.Item(170, 5).Value = "LP1912"
$ws.Cells.Item(171, 1).Value = "10:58:14"
$ws.Cells.Item(171, 2).Value = "11:20"
$ws.Cells.Item(171, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(171, 4).Value = 22
$ws.Cells.Item(171, 5).Value = "LP1912"
$ws.Cells.Item(172, 1).Value = "10:58:14"
$ws.Cells.Item(172, 2).Value = "11:21"
$ws.Cells.Item(172, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(172, 4).Value = 23
$ws.Cells.Item(172, 5).Value = "LP1912"
$ws.Cells.Item(173, 1).Value = "11:27:22"
$ws.Cells.Item(173, 2).Value = "11:27"
$ws.Cells.Item(173, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = "LP1912"
$ws.Cells.Item(174, 1).Value = "11:27:22"
$ws.Cells.Item(174, 2).Value = "11:28"
$ws.Cells.Item(174, 3).Value = "15_ABASTO"
$ws.Cells.Item(174, 4).Value = 1
$ws.Cells.Item(174, 5).Value = "LP1912"
$ws.Cells.Item(175, 1).Value = "11:27:22"
$ws.Cells.Item(175, 2).Value = "11:28"
$ws.Cells.Item(175, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(175, 4).Value = 1
$ws.Cells.Item(175, 5).Value = "LP1912"
$ws.Cells.Item(176, 1).Value = "11:27:22"
$ws.Cells.Item(176, 2).Value = "11:32"
$ws.Cells.Item(176, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(176, 4).Value = 5
$ws.Cells.Item(176, 5).Value = "LP1912"
$ws.Cells.Item(177, 1).Value = "11:27:22"
$ws.Cells.Item(177, 2).Value = "11:34"
$ws.Cells.Item(177, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(177, 4).Value = 7
$ws.Cells.Item(177, 5).Value = "LP1912"
$ws.Cells.Item(178, 1).Value = "10:58:14"
$ws.Cells.Item(178, 2).Value = "11:35"
$ws.Cells.Item(178, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(178, 4).Value = 37
$ws.Cells.Item(178, 5).Value = "LP1912"
$ws.Cells.Item(179, 1).Value = "11:27:22"
$ws.Cells.Item(179, 2).Value = "11:35"
$ws.Cells.Item(179, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(179, 4).Value = 8
$ws.Cells.Item(179, 5).Value = "LP1912"
$ws.Cells.Item(180, 1).Value = "10:58:14"
$ws.Cells.Item(180, 2).Value = "11:36"
$ws.Cells.Item(180, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(180, 4).Value = 38
$ws.Cells.Item(180, 5).Value = "LP1912"
$ws.Cells.Item(181, 1).Value = "11:27:22"
$ws.Cells.Item(181, 2).Value = "11:42"
$ws.Cells.Item(181, 3).Value = "17_ROMERO"
$ws.Cells.Item(181, 4).Value = 15
$ws.Cells.Item(181, 5).Value = "LP1912"
$ws.Cells.Item(182, 1).Value = "11:27:22"
$ws.Cells.Item(182, 2).Value = "11:43"
$ws.Cells.Item(182, 3).Value = "10_OLMOS"
$ws.Cells.Item(182, 4).Value = 16
$ws.Cells.Item(182, 5).Value = "LP1912"
$ws.Cells.Item(183, 1).Value = "11:27:22"
$ws.Cells.Item(183, 2).Value = "11:46"
$ws.Cells.Item(183, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(183, 4).Value = 19
$ws.Cells.Item(183, 5).Value = "LP1912"
$ws.Cells.Item(184, 1).Value = "10:25:24"
$ws.Cells.Item(184, 2).Value = "11:49"
$ws.Cells.Item(184, 3).Value = "10_OLMOS"
$ws.Cells.Item(184, 4).Value = 84
$ws.Cells.Item(184, 5).Value = "LP1912"
$ws.Cells.Item(185, 1).Value = "11:27:22"
$ws.Cells.Item(185, 2).Value = "11:51"
$ws.Cells.Item(185, 3).Value = "215B_EL PATO"
$ws.Cells.Item(185, 4).Value = 24
$ws.Cells.Item(185, 5).Value = "LP1912"
$ws.Cells.Item(186, 1).Value = "10:58:14"
$ws.Cells.Item(186, 2).Value = "11:52"
$ws.Cells.Item(186, 3).Value = "15_ABASTO"
$ws.Cells.Item(186, 4).Value = 54
$ws.Cells.Item(186, 5).Value = "LP1912"
$ws.Cells.Item(187, 1).Value = "11:54:22"
$ws.Cells.Item(187, 2).Value = "11:59"
$ws.Cells.Item(187, 3).Value = "225_GOMEZ"
$ws.Cells.Item(187, 4).Value = 5
$ws.Cells.Item(187, 5).Value = "LP1912"
$ws.Cells.Item(188, 1).Value = "10:25:24"
$ws.Cells.Item(188, 2).Value = "12:01"
$ws.Cells.Item(188, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(188, 4).Value = 96
$ws.Cells.Item(188, 5).Value = "LP1912"
$ws.Cells.Item(189, 1).Value = "11:54:22"
$ws.Cells.Item(189, 2).Value = "12:02"
$ws.Cells.Item(189, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(189, 4).Value = 8
$ws.Cells.Item(189, 5).Value = "LP1912"
$ws.Cells.Item(190, 1).Value = "11:54:22"
$ws.Cells.Item(190, 2).Value = "12:05"
$ws.Cells.Item(190, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(190, 4).Value = 11
$ws.Cells.Item(190, 5).Value = "LP1912"
$ws.Cells.Item(191, 1).Value = "11:54:22"
$ws.Cells.Item(191, 2).Value = "12:06"
$ws.Cells.Item(191, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(191, 4).Value = 12
$ws.Cells.Item(191, 5).Value = "LP1912"
$ws.Cells.Item(192, 1).Value = "11:54:22"
$ws.Cells.Item(192, 2).Value = "12:06"
$ws.Cells.Item(192, 3).Value = "14_ABASTO"
$ws.Cells.Item(192, 4).Value = 12
$ws.Cells.Item(192, 5).Value = "LP1912"
$ws.Cells.Item(193, 1).Value = "10:58:14"
$ws.Cells.Item(193, 2).Value = "12:07"
$ws.Cells.Item(193, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(193, 4).Value = 69
$ws.Cells.Item(193, 5).Value = "LP1912"
$ws.Cells.Item(194, 1).Value = "10:58:14"
$ws.Cells.Item(194, 2).Value = "12:07"
$ws.Cells.Item(194, 3).Value = "14_ABASTO"
$ws.Cells.Item(194, 4).Value = 69
$ws.Cells.Item(194, 5).Value = "LP1912"
$ws.Cells.Item(195, 1).Value = "11:54:22"
$ws.Cells.Item(195, 2).Value = "12:14"
$ws.Cells.Item(195, 3).Value = "17_ROMERO"
$ws.Cells.Item(195, 4).Value = 20
$ws.Cells.Item(195, 5).Value = "LP1912"
$ws.Cells.Item(196, 1).Value = "11:54:22"
$ws.Cells.Item(196, 2).Value = "12:17"
$ws.Cells.Item(196, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(196, 4).Value = 23
$ws.Cells.Item(196, 5).Value = "LP1912"
$ws.Cells.Item(197, 1).Value = "11:27:22"
$ws.Cells.Item(197, 2).Value = "12:20"
$ws.Cells.Item(197, 3).Value = "14_ABASTO"
$ws.Cells.Item(197, 4).Value = 53
$ws.Cells.Item(197, 5).Value = "LP1912"
$ws.Cells.Item(198, 1).Value = "11:54:22"
$ws.Cells.Item(198, 2).Value = "12:20"
$ws.Cells.Item(198, 3).Value = "215A_EL PATO"
$ws.Cells.Item(198, 4).Value = 26
$ws.Cells.Item(198, 5).Value = "LP1912"
$ws.Cells.Item(199, 1).Value = "12:20:14"
$ws.Cells.Item(199, 2).Value = "12:21"
$ws.Cells.Item(199, 3).Value = "14_ABASTO"
$ws.Cells.Item(199, 4).Value = 1
$ws.Cells.Item(199, 5).Value = "LP1912"
$ws.Cells.Item(200, 1).Value = "12:20:14"
$ws.Cells.Item(200, 2).Value = "12:21"
$ws.Cells.Item(200, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(200, 4).Value = 1
$ws.Cells.Item(200, 5).Value = "LP1912"
$ws.Cells.Item(201, 1).Value = "12:20:14"
$ws.Cells.Item(201, 2).Value = "12:21"
$ws.Cells.Item(201, 3).Value = "215A_EL PATO"
$ws.Cells.Item(201, 4).Value = 1
$ws.Cells.Item(201, 5).Value = "LP1912"
$ws.Cells.Item(202, 1).Value = "12:20:14"
$ws.Cells.Item(202, 2).Value = "12:21"
$ws.Cells.Item(202, 3).Value = "17_ROMERO"
$ws.Cells.Item(202, 4).Value = 1
$ws.Cells.Item(202, 5).Value = "LP1912"
$ws.Cells.Item(203, 1).Value = "11:27:22"
$ws.Cells.Item(203, 2).Value = "12:34"
$ws.Cells.Item(203, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(203, 4).Value = 67
$ws.Cells.Item(203, 5).Value = "LP1912"
$ws.Cells.Item(204, 1).Value = "12:20:14"
$ws.Cells.Item(204, 2).Value = "12:35"
$ws.Cells.Item(204, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(204, 4).Value = 15
$ws.Cells.Item(204, 5).Value = "LP1912"
$ws.Cells.Item(205, 1).Value = "12:20:14"
$ws.Cells.Item(205, 2).Value = "12:35"
$ws.Cells.Item(205, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(205, 4).Value = 15
$ws.Cells.Item(205, 5).Value = "LP1912"
$ws.Cells.Item(206, 1).Value = "12:20:14"
$ws.Cells.Item(206, 2).Value = "12:37"
$ws.Cells.Item(206, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(206, 4).Value = 17
$ws.Cells.Item(206, 5).Value = "LP1912"
$ws.Cells.Item(207, 1).Value = "11:54:22"
$ws.Cells.Item(207, 2).Value = "12:38"
$ws.Cells.Item(207, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(207, 4).Value = 44
$ws.Cells.Item(207, 5).Value = "LP1912"
$ws.Cells.Item(208, 1).Value = "12:20:14"
$ws.Cells.Item(208, 2).Value = "12:38"
$ws.Cells.Item(208, 3).Value = "17_179 Y 38"
$ws.Cells.Item(208, 4).Value = 18
$ws.Cells.Item(208, 5).Value = "LP1912"
$ws.Cells.Item(209, 1).Value = "11:54:22"
$ws.Cells.Item(209, 2).Value = "12:41"
$ws.Cells.Item(209, 3).Value = "10_OLMOS"
$ws.Cells.Item(209, 4).Value = 47
$ws.Cells.Item(209, 5).Value = "LP1912"
$ws.Cells.Item(210, 1).Value = "12:20:14"
$ws.Cells.Item(210, 2).Value = "12:44"
$ws.Cells.Item(210, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(210, 4).Value = 24
$ws.Cells.Item(210, 5).Value = "LP1912"
$ws.Cells.Item(211, 1).Value = "11:54:22"
$ws.Cells.Item(211, 2).Value = "12:48"
$ws.Cells.Item(211, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(211, 4).Value = 54
$ws.Cells.Item(211, 5).Value = "LP1912"
$ws.Cells.Item(212, 1).Value = "12:48:04"
$ws.Cells.Item(212, 2).Value = "12:49"
$ws.Cells.Item(212, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(212, 4).Value = 1
$ws.Cells.Item(212, 5).Value = "LP1912"
$ws.Cells.Item(213, 1).Value = "12:48:04"
$ws.Cells.Item(213, 2).Value = "12:49"
$ws.Cells.Item(213, 3).Value = "17_ROMERO"
$ws.Cells.Item(213, 4).Value = 1
$ws.Cells.Item(213, 5).Value = "LP1912"
$ws.Cells.Item(214, 1).Value = "12:48:04"
$ws.Cells.Item(214, 2).Value = "12:55"
$ws.Cells.Item(214, 3).Value = "10_OLMOS"
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = "LP1912"
$ws.Cells.Item(215, 1).Value = "12:48:04"
$ws.Cells.Item(215, 2).Value = "13:02"
$ws.Cells.Item(215, 3).Value = "15_ABASTO"
$ws.Cells.Item(215, 4).Value = 14
$ws.Cells.Item(215, 5).Value = "LP1912"
$ws.Cells.Item(216, 1).Value = "12:48:04"
$ws.Cells.Item(216, 2).Value = "13:03"
$ws.Cells.Item(216, 3).Value = "14_ABASTO"
$ws.Cells.Item(216, 4).Value = 15
$ws.Cells.Item(216, 5).Value = "LP1912"
$ws.Cells.Item(217, 1).Value = "11:54:22"
$ws.Cells.Item(217, 2).Value = "13:04"
$ws.Cells.Item(217, 3).Value = "10_OLMOS"
$ws.Cells.Item(217, 4).Value = 70
$ws.Cells.Item(217, 5).Value = "LP1912"
$ws.Cells.Item(218, 1).Value = "13:02:04"
$ws.Cells.Item(218, 2).Value = "13:05"
$ws.Cells.Item(218, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(218, 4).Value = 3
$ws.Cells.Item(218, 5).Value = "LP1912"
$ws.Cells.Item(219, 1).Value = "11:54:22"
$ws.Cells.Item(219, 2).Value = "13:06"
$ws.Cells.Item(219, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(219, 4).Value = 72
$ws.Cells.Item(219, 5).Value = "LP1912"
$ws.Cells.Item(220, 1).Value = "13:02:04"
$ws.Cells.Item(220, 2).Value = "13:07"
$ws.Cells.Item(220, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(220, 4).Value = 5
$ws.Cells.Item(220, 5).Value = "LP1912"
$ws.Cells.Item(221, 1).Value = "13:02:04"
$ws.Cells.Item(221, 2).Value = "13:08"
$ws.Cells.Item(221, 3).Value = "10_OLMOS"
$ws.Cells.Item(221, 4).Value = 6
$ws.Cells.Item(221, 5).Value = "LP1912"
$ws.Cells.Item(222, 1).Value = "12:48:04"
$ws.Cells.Item(222, 2).Value = "13:10"
$ws.Cells.Item(222, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(222, 4).Value = 22
$ws.Cells.Item(222, 5).Value = "LP1912"
$ws.Cells.Item(223, 1).Value = "13:02:04"
$ws.Cells.Item(223, 2).Value = "13:14"
$ws.Cells.Item(223, 3).Value = "215D_EL PATO"
$ws.Cells.Item(223, 4).Value = 12
$ws.Cells.Item(223, 5).Value = "LP1912"
$ws.Cells.Item(224, 1).Value = "13:02:04"
$ws.Cells.Item(224, 2).Value = "13:15"
$ws.Cells.Item(224, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(224, 4).Value = 13
$ws.Cells.Item(224, 5).Value = "LP1912"
$ws.Cells.Item(225, 1).Value = "12:20:14"
$ws.Cells.Item(225, 2).Value = "13:19"
$ws.Cells.Item(225, 3).Value = "10_OLMOS"
$ws.Cells.Item(225, 4).Value = 59
$ws.Cells.Item(225, 5).Value = "LP1912"
$ws.Cells.Item(226, 1).Value = "13:02:04"
$ws.Cells.Item(226, 2).Value = "13:20"
$ws.Cells.Item(226, 3).Value = "10_OLMOS"
$ws.Cells.Item(226, 4).Value = 18
$ws.Cells.Item(226, 5).Value = "LP1912"
$ws.Cells.Item(227, 1).Value = "13:02:04"
$ws.Cells.Item(227, 2).Value = "13:21"
$ws.Cells.Item(227, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(227, 4).Value = 19
$ws.Cells.Item(227, 5).Value = "LP1912"
$ws.Cells.Item(228, 1).Value = "13:02:04"
$ws.Cells.Item(228, 2).Value = "13:26"
$ws.Cells.Item(228, 3).Value = "15_ABASTO"
$ws.Cells.Item(228, 4).Value = 24
$ws.Cells.Item(228, 5).Value = "LP1912"
$ws.Cells.Item(229, 1).Value = "11:54:22"
$ws.Cells.Item(229, 2).Value = "13:26"
$ws.Cells.Item(229, 3).Value = "14_ABASTO"
$ws.Cells.Item(229, 4).Value = 92
$ws.Cells.Item(229, 5).Value = "LP1912"
$ws.Cells.Item(230, 1).Value = "13:02:04"
$ws.Cells.Item(230, 2).Value = "13:27"
$ws.Cells.Item(230, 3).Value = "14_ABASTO"
$ws.Cells.Item(230, 4).Value = 25
$ws.Cells.Item(230, 5).Value = "LP1912"
$ws.Cells.Item(231, 1).Value = "13:02:04"
$ws.Cells.Item(231, 2).Value = "13:34"
$ws.Cells.Item(231, 3).Value = "10_OLMOS"
$ws.Cells.Item(231, 4).Value = 32
$ws.Cells.Item(231, 5).Value = "LP1912"
$ws.Cells.Item(232, 1).Value = "13:02:04"
$ws.Cells.Item(232, 2).Value = "13:35"
$ws.Cells.Item(232, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(232, 4).Value = 33
$ws.Cells.Item(232, 5).Value = "LP1912"
$ws.Cells.Item(233, 1).Value = "13:02:04"
$ws.Cells.Item(233, 2).Value = "13:36"
$ws.Cells.Item(233, 3).Value = "15_ABASTO"
$ws.Cells.Item(233, 4).Value = 34
$ws.Cells.Item(233, 5).Value = "LP1912"
$ws.Cells.Item(234, 1).Value = "13:02:04"
$ws.Cells.Item(234, 2).Value = "13:46"
$ws.Cells.Item(234, 3).Value = "17_ROMERO"
$ws.Cells.Item(234, 4).Value = 44
$ws.Cells.Item(234, 5).Value = "LP1912"
$ws.Cells.Item(235, 1).Value = "13:02:04"
$ws.Cells.Item(235, 2).Value = "13:47"
$ws.Cells.Item(235, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(235, 4).Value = 45
$ws.Cells.Item(235, 5).Value = "LP1912"
$ws.Cells.Item(236, 1).Value = "13:49:12"
$ws.Cells.Item(236, 2).Value = "13:49"
$ws.Cells.Item(236, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(236, 4).Value = 0
$ws.Cells.Item(236, 5).Value = "LP1912"
$ws.Cells.Item(237, 1).Value = "13:49:12"
$ws.Cells.Item(237, 2).Value = "13:50"
$ws.Cells.Item(237, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(237, 4).Value = 1
$ws.Cells.Item(237, 5).Value = "LP1912"
$ws.Cells.Item(238, 1).Value = "11:54:22"
$ws.Cells.Item(238, 2).Value = "13:50"
$ws.Cells.Item(238, 3).Value = "215A_EL PATO"
$ws.Cells.Item(238, 4).Value = 116
$ws.Cells.Item(238, 5).Value = "LP1912"
$ws.Cells.Item(239, 1).Value = "13:02:04"
$ws.Cells.Item(239, 2).Value = "13:51"
$ws.Cells.Item(239, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(239, 4).Value = 49
$ws.Cells.Item(239, 5).Value = "LP1912"
$ws.Cells.Item(240, 1).Value = "13:49:12"
$ws.Cells.Item(240, 2).Value = "13:51"
$ws.Cells.Item(240, 3).Value = "215A_EL PATO"
$ws.Cells.Item(240, 4).Value = 2
$ws.Cells.Item(240, 5).Value = "LP1912"
$ws.Cells.Item(241, 1).Value = "13:49:12"
$ws.Cells.Item(241, 2).Value = "13:56"
$ws.Cells.Item(241, 3).Value = "225_GOMEZ"
$ws.Cells.Item(241, 4).Value = 7
$ws.Cells.Item(241, 5).Value = "LP1912"
$ws.Cells.Item(242, 1).Value = "13:49:12"
$ws.Cells.Item(242, 2).Value = "13:57"
$ws.Cells.Item(242, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(242, 4).Value = 8
$ws.Cells.Item(242, 5).Value = "LP1912"
$ws.Cells.Item(243, 1).Value = "13:49:12"
$ws.Cells.Item(243, 2).Value = "14:04"
$ws.Cells.Item(243, 3).Value = "17_ROMERO"
$ws.Cells.Item(243, 4).Value = 15
$ws.Cells.Item(243, 5).Value = "LP1912"
$ws.Cells.Item(244, 1).Value = "13:49:12"
$ws.Cells.Item(244, 2).Value = "14:05"
$ws.Cells.Item(244, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(244, 4).Value = 16
$ws.Cells.Item(244, 5).Value = "LP1912"
$ws.Cells.Item(245, 1).Value = "13:49:12"
$ws.Cells.Item(245, 2).Value = "14:09"
$ws.Cells.Item(245, 3).Value = "10_OLMOS"
$ws.Cells.Item(245, 4).Value = 20
$ws.Cells.Item(245, 5).Value = "LP1912"
$ws.Cells.Item(246, 1).Value = "13:49:12"
$ws.Cells.Item(246, 2).Value = "14:12"
$ws.Cells.Item(246, 3).Value = "15_ABASTO"
$ws.Cells.Item(246, 4).Value = 23
$ws.Cells.Item(246, 5).Value = "LP1912"
$ws.Cells.Item(247, 1).Value = "14:17:09"
$ws.Cells.Item(247, 2).Value = "14:17"
$ws.Cells.Item(247, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(247, 4).Value = 0
$ws.Cells.Item(247, 5).Value = "LP1912"
$ws.Cells.Item(248, 1).Value = "14:17:09"
$ws.Cells.Item(248, 2).Value = "14:20"
$ws.Cells.Item(248, 3).Value = "215C_EL PATO"
$ws.Cells.Item(248, 4).Value = 3
$ws.Cells.Item(248, 5).Value = "LP1912"
$ws.Cells.Item(249, 1).Value = "14:17:09"
$ws.Cells.Item(249, 2).Value = "14:21"
$ws.Cells.Item(249, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(249, 4).Value = 4
$ws.Cells.Item(249, 5).Value = "LP1912"
$ws.Cells.Item(250, 1).Value = "14:17:09"
$ws.Cells.Item(250, 2).Value = "14:27"
$ws.Cells.Item(250, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(250, 4).Value = 10
$ws.Cells.Item(250, 5).Value = "LP1912"
$ws.Cells.Item(251, 1).Value = "14:17:09"
$ws.Cells.Item(251, 2).Value = "14:28"
$ws.Cells.Item(251, 3).Value = "15_ABASTO"
$ws.Cells.Item(251, 4).Value = 11
$ws.Cells.Item(251, 5).Value = "LP1912"
$ws.Cells.Item(252, 1).Value = "14:17:09"
$ws.Cells.Item(252, 2).Value = "14:35"
$ws.Cells.Item(252, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(252, 4).Value = 18
$ws.Cells.Item(252, 5).Value = "LP1912"
$ws.Cells.Item(253, 1).Value = "14:17:09"
$ws.Cells.Item(253, 2).Value = "14:44"
$ws.Cells.Item(253, 3).Value = "15_ABASTO"
$ws.Cells.Item(253, 4).Value = 27
$ws.Cells.Item(253, 5).Value = "LP1912"
$ws.Cells.Item(254, 1).Value = "14:17:09"
$ws.Cells.Item(254, 2).Value = "14:45"
$ws.Cells.Item(254, 3).Value = "14_ABASTO"
$ws.Cells.Item(254, 4).Value = 28
$ws.Cells.Item(254, 5).Value = "LP1912"
$ws.Cells.Item(255, 1).Value = "14:47:31"
$ws.Cells.Item(255, 2).Value = "14:47"
$ws.Cells.Item(255, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(255, 4).Value = 0
$ws.Cells.Item(255, 5).Value = "LP1912"
$ws.Cells.Item(256, 1).Value = "14:47:31"
$ws.Cells.Item(256, 2).Value = "14:48"
$ws.Cells.Item(256, 3).Value = "15_ABASTO"
$ws.Cells.Item(256, 4).Value = 1
$ws.Cells.Item(256, 5).Value = "LP1912"
$ws.Cells.Item(257, 1).Value = "13:02:04"
$ws.Cells.Item(257, 2).Value = "14:52"
$ws.Cells.Item(257, 3).Value = "14_ABASTO"
$ws.Cells.Item(257, 4).Value = 110
$ws.Cells.Item(257, 5).Value = "LP1912"
$ws.Cells.Item(258, 1).Value = "14:47:31"
$ws.Cells.Item(258, 2).Value = "14:56"
$ws.Cells.Item(258, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(258, 4).Value = 9
$ws.Cells.Item(258, 5).Value = "LP1912"
$ws.Cells.Item(259, 1).Value = "14:17:09"
$ws.Cells.Item(259, 2).Value = "14:57"
$ws.Cells.Item(259, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(259, 4).Value = 40
$ws.Cells.Item(259, 5).Value = "LP1912"
$ws.Cells.Item(260, 1).Value = "14:47:31"
$ws.Cells.Item(260, 2).Value = "14:58"
$ws.Cells.Item(260, 3).Value = "215B_EL PATO"
$ws.Cells.Item(260, 4).Value = 11
$ws.Cells.Item(260, 5).Value = "LP1912"
$ws.Cells.Item(261, 1).Value = "14:47:31"
$ws.Cells.Item(261, 2).Value = "15:00"
$ws.Cells.Item(261, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(261, 4).Value = 13
$ws.Cells.Item(261, 5).Value = "LP1912"
$ws.Cells.Item(262, 1).Value = "14:47:31"
$ws.Cells.Item(262, 2).Value = "15:04"
$ws.Cells.Item(262, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(262, 4).Value = 17
$ws.Cells.Item(262, 5).Value = "LP1912"
$ws.Cells.Item(263, 1).Value = "14:47:31"
$ws.Cells.Item(263, 2).Value = "15:05"
$ws.Cells.Item(263, 3).Value = "10_OLMOS"
$ws.Cells.Item(263, 4).Value = 18
$ws.Cells.Item(263, 5).Value = "LP1912"
$ws.Cells.Item(264, 1).Value = "14:47:31"
$ws.Cells.Item(264, 2).Value = "15:06"
$ws.Cells.Item(264, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(264, 4).Value = 19
$ws.Cells.Item(264, 5).Value = "LP1912"
$ws.Cells.Item(265, 1).Value = "14:47:31"
$ws.Cells.Item(265, 2).Value = "15:10"
$ws.Cells.Item(265, 3).Value = "17_ROMERO"
$ws.Cells.Item(265, 4).Value = 23
$ws.Cells.Item(265, 5).Value = "LP1912"
$ws.Cells.Item(266, 1).Value = "14:47:31"
$ws.Cells.Item(266, 2).Value = "15:13"
$ws.Cells.Item(266, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(266, 4).Value = 26
$ws.Cells.Item(266, 5).Value = "LP1912"
$ws.Cells.Item(267, 1).Value = "14:17:09"
$ws.Cells.Item(267, 2).Value = "15:14"
$ws.Cells.Item(267, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(267, 4).Value = 57
$ws.Cells.Item(267, 5).Value = "LP1912"
$ws.Cells.Item(268, 1).Value = "13:49:12"
$ws.Cells.Item(268, 2).Value = "15:20"
$ws.Cells.Item(268, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(268, 4).Value = 91
$ws.Cells.Item(268, 5).Value = "LP1912"
$ws.Cells.Item(269, 1).Value = "14:47:31"
$ws.Cells.Item(269, 2).Value = "15:21"
$ws.Cells.Item(269, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(269, 4).Value = 34
$ws.Cells.Item(269, 5).Value = "LP1912"
$ws.Cells.Item(270, 1).Value = "14:47:31"
$ws.Cells.Item(270, 2).Value = "15:32"
$ws.Cells.Item(270, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(270, 4).Value = 45
$ws.Cells.Item(270, 5).Value = "LP1912"
$ws.Cells.Item(271, 1).Value = "13:49:12"
$ws.Cells.Item(271, 2).Value = "15:35"
$ws.Cells.Item(271, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(271, 4).Value = 106
$ws.Cells.Item(271, 5).Value = "LP1912"
$ws.Cells.Item(272, 1).Value = "14:47:31"
$ws.Cells.Item(272, 2).Value = "15:37"
$ws.Cells.Item(272, 3).Value = "10_OLMOS"
$ws.Cells.Item(272, 4).Value = 50
$ws.Cells.Item(272, 5).Value = "LP1912"
$ws.Cells.Item(273, 1).Value = "14:47:31"
$ws.Cells.Item(273, 2).Value = "15:38"
$ws.Cells.Item(273, 3).Value = "215A_EL PATO"
$ws.Cells.Item(273, 4).Value = 51
$ws.Cells.Item(273, 5).Value = "LP1912"
$ws.Cells.Item(274, 1).Value = "14:17:09"
$ws.Cells.Item(274, 2).Value = "15:38"
$ws.Cells.Item(274, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(274, 4).Value = 81
$ws.Cells.Item(274, 5).Value = "LP1912"
$ws.Cells.Item(275, 1).Value = "14:17:09"
$ws.Cells.Item(275, 2).Value = "15:39"
$ws.Cells.Item(275, 3).Value = "215A_EL PATO"
$ws.Cells.Item(275, 4).Value = 82
$ws.Cells.Item(275, 5).Value = "LP1912"
$ws.Cells.Item(276, 1).Value = "14:17:09"
$ws.Cells.Item(276, 2).Value = "15:40"
$ws.Cells.Item(276, 3).Value = "14_ABASTO"
$ws.Cells.Item(276, 4).Value = 83
$ws.Cells.Item(276, 5).Value = "LP1912"
$ws.Cells.Item(277, 1).Value = "14:47:31"
$ws.Cells.Item(277, 2).Value = "15:42"
$ws.Cells.Item(277, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(277, 4).Value = 55
$ws.Cells.Item(277, 5).Value = "LP1912"
$ws.Cells.Item(278, 1).Value = "14:47:31"
$ws.Cells.Item(278, 2).Value = "15:45"
$ws.Cells.Item(278, 3).Value = "14_ABASTO"
$ws.Cells.Item(278, 4).Value = 58
$ws.Cells.Item(278, 5).Value = "LP1912"
$ws.Cells.Item(279, 1).Value = "14:47:31"
$ws.Cells.Item(279, 2).Value = "15:46"
$ws.Cells.Item(279, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(279, 4).Value = 59
$ws.Cells.Item(279, 5).Value = "LP1912"
$ws.Cells.Item(280, 1).Value = "14:17:09"
$ws.Cells.Item(280, 2).Value = "15:47"
$ws.Cells.Item(280, 3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(280, 4).Value = 90
$ws.Cells.Item(280, 5).Value = "LP1912"
$ws.Cells.Item(281, 1).Value = "14:47:31"
$ws.Cells.Item(281, 2).Value = "15:52"
$ws.Cells.Item(281, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(281, 4).Value = 65
$ws.Cells.Item(281, 5).Value = "LP1912"
$ws.Cells.Item(282, 1).Value = "14:47:31"
$ws.Cells.Item(282, 2).Value = "15:53"
$ws.Cells.Item(282, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(282, 4).Value = 66
$ws.Cells.Item(282, 5).Value = "LP1912"
$ws.Cells.Item(283, 1).Value = "14:17:09"
$ws.Cells.Item(283, 2).Value = "15:54"
$ws.Cells.Item(283, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(283, 4).Value = 97
$ws.Cells.Item(283, 5).Value = "LP1912"
$ws.Cells.Item(284, 1).Value = "14:47:31"
$ws.Cells.Item(284, 2).Value = "15:56"
$ws.Cells.Item(284, 3).Value = "17_ROMERO"
$ws.Cells.Item(284, 4).Value = 69
$ws.Cells.Item(284, 5).Value = "LP1912"
$ws.Cells.Item(285, 1).Value = "14:17:09"
$ws.Cells.Item(285, 2).Value = "15:57"
$ws.Cells.Item(285, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(285, 4).Value = 100
$ws.Cells.Item(285, 5).Value = "LP1912"
$ws.Cells.Item(286, 1).Value = "14:47:31"
$ws.Cells.Item(286, 2).Value = "16:08"
$ws.Cells.Item(286, 3).Value = "14_ABASTO"
$ws.Cells.Item(286, 4).Value = 81
$ws.Cells.Item(286, 5).Value = "LP1912"
$ws.Cells.Item(287, 1).Value = "14:47:31"
$ws.Cells.Item(287, 2).Value = "16:15"
$ws.Cells.Item(287, 3).Value = "225_C ROCA-H SUR"
$ws.Cells.Item(287, 4).Value = 88
$ws.Cells.Item(287, 5).Value = "LP1912"
$ws.Cells.Item(288, 1).Value = "14:47:31"
$ws.Cells.Item(288, 2).Value = "16:20"
$ws.Cells.Item(288, 3).Value = "215C_EL PATO"
$ws.Cells.Item(288, 4).Value = 93
$ws.Cells.Item(288, 5).Value = "LP1912"
$ws.Cells.Item(289, 1).Value = "14:47:31"
$ws.Cells.Item(289, 2).Value = "16:21"
$ws.Cells.Item(289, 3).Value = "26_HERNANDEZ"
$ws.Cells.Item(289, 4).Value = 94
$ws.Cells.Item(289, 5).Value = "LP1912"
$ws.Cells.Item(290, 1).Value = "14:47:31"
$ws.Cells.Item(290, 2).Value = "16:30"
$ws.Cells.Item(290, 3).Value = "15_ABASTO"
$ws.Cells.Item(290, 4).Value = 103
$ws.Cells.Item(290, 5).Value = "LP1912"
$ws.Cells.Item(291, 1).Value = "14:47:31"
$ws.Cells.Item(291, 2).Value = "16:43"
$ws.Cells.Item(291, 3).Value = "225_GOMEZ"
$ws.Cells.Item(291, 4).Value = 116
$ws.Cells.Item(291, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(1, 1).Value = "LÍNEA 141 - LP1912-215 - 26/01/2026"
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:47:32"
$ws.Cells.Item(3, 1).Value = "Total filas: 32"
$ws.Cells.Item(5, 1).Value = "Hora_Scrap"
$ws.Cells.Item(5, 2).Value = "Hora_Llegada"
$ws.Cells.Item(5, 3).Value = "Linea"
$ws.Cells.Item(5, 4).Value = "Minutos"
$ws.Cells.Item(5, 5).Value = "Parada"
$ws.Cells.Item(6, 1).Value = "04:45:48"
$ws.Cells.Item(6, 2).Value = "04:45"
$ws.Cells.Item(6, 3).Value = "215A_EL PATO"
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = "LP1912"
$ws.Cells.Item(7, 1).Value = "04:21:09"
$ws.Cells.Item(7, 2).Value = "04:46"
$ws.Cells.Item(7, 3).Value = "215A_EL PATO"
$ws.Cells.Item(7, 4).Value = 25
$ws.Cells.Item(7, 5).Value = "LP1912"
$ws.Cells.Item(8, 1).Value = "04:56:30"
$ws.Cells.Item(8, 2).Value = "05:34"
$ws.Cells.Item(8, 3).Value = "215B_EL PATO"
$ws.Cells.Item(8, 4).Value = 38
$ws.Cells.Item(8, 5).Value = "LP1912"
$ws.Cells.Item(9, 1).Value = "05:24:16"
$ws.Cells.Item(9, 2).Value = "05:35"
$ws.Cells.Item(9, 3).Value = "215B_EL PATO"
$ws.Cells.Item(9, 4).Value = 11
$ws.Cells.Item(9, 5).Value = "LP1912"
$ws.Cells.Item(10, 1).Value = "05:24:16"
$ws.Cells.Item(10, 2).Value = "06:11"
$ws.Cells.Item(10, 3).Value = "215A_EL PATO"
$ws.Cells.Item(10, 4).Value = 47
$ws.Cells.Item(10, 5).Value = "LP1912"
$ws.Cells.Item(11, 1).Value = "05:55:02"
$ws.Cells.Item(11, 2).Value = "06:12"
$ws.Cells.Item(11, 3).Value = "215A_EL PATO"
$ws.Cells.Item(11, 4).Value = 17
$ws.Cells.Item(11, 5).Value = "LP1912"
$ws.Cells.Item(12, 1).Value = "06:25:28"
$ws.Cells.Item(12, 2).Value = "06:46"
$ws.Cells.Item(12, 3).Value = "215C_EL PATO"
$ws.Cells.Item(12, 4).Value = 21
$ws.Cells.Item(12, 5).Value = "LP1912"
$ws.Cells.Item(13, 1).Value = "05:55:02"
$ws.Cells.Item(13, 2).Value = "06:47"
$ws.Cells.Item(13, 3).Value = "215C_EL PATO"
$ws.Cells.Item(13, 4).Value = 52
$ws.Cells.Item(13, 5).Value = "LP1912"
$ws.Cells.Item(14, 1).Value = "06:54:06"
$ws.Cells.Item(14, 2).Value = "06:55"
$ws.Cells.Item(14, 3).Value = "215C_EL PATO"
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = "LP1912"
$ws.Cells.Item(15, 1).Value = "06:25:28"
$ws.Cells.Item(15, 2).Value = "07:11"
$ws.Cells.Item(15, 3).Value = "215A_EL PATO"
$ws.Cells.Item(15, 4).Value = 46
$ws.Cells.Item(15, 5).Value = "LP1912"
$ws.Cells.Item(16, 1).Value = "06:54:06"
$ws.Cells.Item(16, 2).Value = "07:12"
$ws.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws.Cells.Item(16, 4).Value = 18
$ws.Cells.Item(16, 5).Value = "LP1912"
$ws.Cells.Item(17, 1).Value = "07:17:59"
$ws.Cells.Item(17, 2).Value = "07:51"
$ws.Cells.Item(17, 3).Value = "215D_EL PATO"
$ws.Cells.Item(17, 4).Value = 34
$ws.Cells.Item(17, 5).Value = "LP1912"
$ws.Cells.Item(18, 1).Value = "07:48:05"
$ws.Cells.Item(18, 2).Value = "07:52"
$ws.Cells.Item(18, 3).Value = "215D_EL PATO"
$ws.Cells.Item(18, 4).Value = 4
$ws.Cells.Item(18, 5).Value = "LP1912"
$ws.Cells.Item(19, 1).Value = "07:17:59"
$ws.Cells.Item(19, 2).Value = "08:22"
$ws.Cells.Item(19, 3).Value = "215B_EL PATO"
$ws.Cells.Item(19, 4).Value = 65
$ws.Cells.Item(19, 5).Value = "LP1912"
$ws.Cells.Item(20, 1).Value = "08:01:10"
$ws.Cells.Item(20, 2).Value = "08:23"
$ws.Cells.Item(20, 3).Value = "215B_EL PATO"
$ws.Cells.Item(20, 4).Value = 22
$ws.Cells.Item(20, 5).Value = "LP1912"
$ws.Cells.Item(21, 1).Value = "08:47:26"
$ws.Cells.Item(21, 2).Value = "09:01"
$ws.Cells.Item(21, 3).Value = "215A_EL PATO"
$ws.Cells.Item(21, 4).Value = 14
$ws.Cells.Item(21, 5).Value = "LP1912"
$ws.Cells.Item(22, 1).Value = "08:55:01"
$ws.Cells.Item(22, 2).Value = "09:02"
$ws.Cells.Item(22, 3).Value = "215A_EL PATO"
$ws.Cells.Item(22, 4).Value = 7
$ws.Cells.Item(22, 5).Value = "LP1912"
$ws.Cells.Item(23, 1).Value = "09:28:08"
$ws.Cells.Item(23, 2).Value = "09:42"
$ws.Cells.Item(23, 3).Value = "215C_EL PATO"
$ws.Cells.Item(23, 4).Value = 14
$ws.Cells.Item(23, 5).Value = "LP1912"
$ws.Cells.Item(24, 1).Value = "10:25:24"
$ws.Cells.Item(24, 2).Value = "10:26"
$ws.Cells.Item(24, 3).Value = "215A_EL PATO"
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = "LP1912"
$ws.Cells.Item(25, 1).Value = "09:28:08"
$ws.Cells.Item(25, 2).Value = "10:27"
$ws.Cells.Item(25, 3).Value = "215A_EL PATO"
$ws.Cells.Item(25, 4).Value = 59
$ws.Cells.Item(25, 5).Value = "LP1912"
$ws.Cells.Item(26, 1).Value = "10:58:14"
$ws.Cells.Item(26, 2).Value = "11:02"
$ws.Cells.Item(26, 3).Value = "215C_EL PATO"
$ws.Cells.Item(26, 4).Value = 4
$ws.Cells.Item(26, 5).Value = "LP1912"
$ws.Cells.Item(27, 1).Value = "11:27:22"
$ws.Cells.Item(27, 2).Value = "11:51"
$ws.Cells.Item(27, 3).Value = "215B_EL PATO"
$ws.Cells.Item(27, 4).Value = 24
$ws.Cells.Item(27, 5).Value = "LP1912"
$ws.Cells.Item(28, 1).Value = "11:54:22"
$ws.Cells.Item(28, 2).Value = "12:20"
$ws.Cells.Item(28, 3).Value = "215A_EL PATO"
$ws.Cells.Item(28, 4).Value = 26
$ws.Cells.Item(28, 5).Value = "LP1912"
$ws.Cells.Item(29, 1).Value = "12:20:14"
$ws.Cells.Item(29, 2).Value = "12:21"
$ws.Cells.Item(29, 3).Value = "215A_EL PATO"
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 5).Value = "LP1912"
$ws.Cells.Item(30, 1).Value = "13:02:04"
$ws.Cells.Item(30, 2).Value = "13:14"
$ws.Cells.Item(30, 3).Value = "215D_EL PATO"
$ws.Cells.Item(30, 4).Value = 12
$ws.Cells.Item(30, 5).Value = "LP1912"
$ws.Cells.Item(31, 1).Value = "11:54:22"
$ws.Cells.Item(31, 2).Value = "13:50"
$ws.Cells.Item(31, 3).Value = "215A_EL PATO"
$ws.Cells.Item(31, 4).Value = 116
$ws.Cells.Item(31, 5).Value = "LP1912"
$ws.Cells.Item(32, 1).Value = "13:49:12"
$ws.Cells.Item(32, 2).Value = "13:51"
$ws.Cells.Item(32, 3).Value = "215A_EL PATO"
$ws.Cells.Item(32, 4).Value = 2
$ws.Cells.Item(32, 5).Value = "LP1912"
$ws.Cells.Item(33, 1).Value = "14:17:09"
$ws.Cells.Item(33, 2).Value = "14:20"
$ws.Cells.Item(33, 3).Value = "215C_EL PATO"
$ws.Cells.Item(33, 4).Value = 3
$ws.Cells.Item(33, 5).Value = "LP1912"
$ws.Cells.Item(34, 1).Value = "14:47:31"
$ws.Cells.Item(34, 2).Value = "14:58"
$ws.Cells.Item(34, 3).Value = "215B_EL PATO"
$ws.Cells.Item(34, 4).Value = 11
$ws.Cells.Item(34, 5).Value = "LP1912"
$ws.Cells.Item(35, 1).Value = "14:47:31"
$ws.Cells.Item(35, 2).Value = "15:38"
$ws.Cells.Item(35, 3).Value = "215A_EL PATO"
$ws.Cells.Item(35, 4).Value = 51
$ws.Cells.Item(35, 5).Value = "LP1912"
$ws.Cells.Item(36, 1).Value = "14:17:09"
$ws.Cells.Item(36, 2).Value = "15:39"
$ws.Cells.Item(36, 3).Value = "215A_EL PATO"
$ws.Cells.Item(36, 4).Value = 82
$ws.Cells.Item(36, 5).Value = "LP1912"
$ws.Cells.Item(37, 1).Value = "14:47:31"
$ws.Cells.Item(37, 2).Value = "16:20"
$ws.Cells.Item(37, 3).Value = "215C_EL PATO"
$ws.Cells.Item(37, 4).Value = 93
$ws.Cells.Item(37, 5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(1, 1).Value = "LÍNEA 141 - 6203-6173 - 26/01/2026"
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:47:32"
$ws.Cells.Item(3, 1).Value = "Total filas: 32"
$ws.Cells.Item(5, 1).Value = "Hora_Scrap"
$ws.Cells.Item(5, 2).Value = "Hora_Llegada"
$ws.Cells.Item(5, 3).Value = "Linea"
$ws.Cells.Item(5, 4).Value = "Minutos"
$ws.Cells.Item(5, 5).Value = "Parada"
$ws.Cells.Item(6, 1).Value = "04:45:48"
$ws.Cells.Item(6, 2).Value = "05:43"
$ws.Cells.Item(6, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(6, 4).Value = 58
$ws.Cells.Item(6, 5).Value = "L6173"
$ws.Cells.Item(7, 1).Value = "05:24:16"
$ws.Cells.Item(7, 2).Value = "05:44"
$ws.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(7, 4).Value = 20
$ws.Cells.Item(7, 5).Value = "L6173"
$ws.Cells.Item(8, 1).Value = "04:45:48"
$ws.Cells.Item(8, 2).Value = "06:08"
$ws.Cells.Item(8, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(8, 4).Value = 83
$ws.Cells.Item(8, 5).Value = "L6173"
$ws.Cells.Item(9, 1).Value = "05:55:02"
$ws.Cells.Item(9, 2).Value = "06:09"
$ws.Cells.Item(9, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(9, 4).Value = 14
$ws.Cells.Item(9, 5).Value = "L6173"
$ws.Cells.Item(10, 1).Value = "04:45:48"
$ws.Cells.Item(10, 2).Value = "06:32"
$ws.Cells.Item(10, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(10, 4).Value = 107
$ws.Cells.Item(10, 5).Value = "L6203"
$ws.Cells.Item(11, 1).Value = "06:25:28"
$ws.Cells.Item(11, 2).Value = "06:33"
$ws.Cells.Item(11, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(11, 4).Value = 8
$ws.Cells.Item(11, 5).Value = "L6203"
$ws.Cells.Item(12, 1).Value = "06:54:06"
$ws.Cells.Item(12, 2).Value = "07:00"
$ws.Cells.Item(12, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(12, 4).Value = 6
$ws.Cells.Item(12, 5).Value = "L6173"
$ws.Cells.Item(13, 1).Value = "07:17:59"
$ws.Cells.Item(13, 2).Value = "07:34"
$ws.Cells.Item(13, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(13, 4).Value = 17
$ws.Cells.Item(13, 5).Value = "L6173"
$ws.Cells.Item(14, 1).Value = "06:54:06"
$ws.Cells.Item(14, 2).Value = "07:35"
$ws.Cells.Item(14, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(14, 4).Value = 41
$ws.Cells.Item(14, 5).Value = "L6173"
$ws.Cells.Item(15, 1).Value = "06:25:28"
$ws.Cells.Item(15, 2).Value = "08:07"
$ws.Cells.Item(15, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(15, 4).Value = 102
$ws.Cells.Item(15, 5).Value = "L6203"
$ws.Cells.Item(16, 1).Value = "06:54:06"
$ws.Cells.Item(16, 2).Value = "08:13"
$ws.Cells.Item(16, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(16, 4).Value = 79
$ws.Cells.Item(16, 5).Value = "L6203"
$ws.Cells.Item(17, 1).Value = "07:17:59"
$ws.Cells.Item(17, 2).Value = "08:19"
$ws.Cells.Item(17, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(17, 4).Value = 62
$ws.Cells.Item(17, 5).Value = "L6203"
$ws.Cells.Item(18, 1).Value = "08:01:10"
$ws.Cells.Item(18, 2).Value = "08:22"
$ws.Cells.Item(18, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(18, 4).Value = 21
$ws.Cells.Item(18, 5).Value = "L6203"
$ws.Cells.Item(19, 1).Value = "07:48:05"
$ws.Cells.Item(19, 2).Value = "08:25"
$ws.Cells.Item(19, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(19, 4).Value = 37
$ws.Cells.Item(19, 5).Value = "L6203"
$ws.Cells.Item(20, 1).Value = "08:31:01"
$ws.Cells.Item(20, 2).Value = "08:33"
$ws.Cells.Item(20, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(20, 4).Value = 2
$ws.Cells.Item(20, 5).Value = "L6203"
$ws.Cells.Item(21, 1).Value = "07:17:59"
$ws.Cells.Item(21, 2).Value = "08:34"
$ws.Cells.Item(21, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(21, 4).Value = 77
$ws.Cells.Item(21, 5).Value = "L6173"
$ws.Cells.Item(22, 1).Value = "08:01:10"
$ws.Cells.Item(22, 2).Value = "08:35"
$ws.Cells.Item(22, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(22, 4).Value = 34
$ws.Cells.Item(22, 5).Value = "L6173"
$ws.Cells.Item(23, 1).Value = "08:31:01"
$ws.Cells.Item(23, 2).Value = "08:36"
$ws.Cells.Item(23, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(23, 4).Value = 5
$ws.Cells.Item(23, 5).Value = "L6173"
$ws.Cells.Item(24, 1).Value = "07:17:59"
$ws.Cells.Item(24, 2).Value = "09:08"
$ws.Cells.Item(24, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(24, 4).Value = 111
$ws.Cells.Item(24, 5).Value = "L6203"
$ws.Cells.Item(25, 1).Value = "08:55:01"
$ws.Cells.Item(25, 2).Value = "09:09"
$ws.Cells.Item(25, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(25, 4).Value = 14
$ws.Cells.Item(25, 5).Value = "L6203"
$ws.Cells.Item(26, 1).Value = "09:28:08"
$ws.Cells.Item(26, 2).Value = "10:03"
$ws.Cells.Item(26, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws.Cells.Item(26, 4).Value = 35
$ws.Cells.Item(26, 5).Value = "L6173"
$ws.Cells.Item(27, 1).Value = "10:25:24"
$ws.Cells.Item(27, 2).Value = "10:54"
$ws.Cells.Item(27, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(27, 4).Value = 29
$ws.Cells.Item(27, 5).Value = "L6173"
$ws.Cells.Item(28, 1).Value = "10:58:14"
$ws.Cells.Item(28, 2).Value = "11:14"
$ws.Cells.Item(28, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(28, 4).Value = 16
$ws.Cells.Item(28, 5).Value = "L6203"
$ws.Cells.Item(29, 1).Value = "11:54:22"
$ws.Cells.Item(29, 2).Value = "12:04"
$ws.Cells.Item(29, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(29, 4).Value = 10
$ws.Cells.Item(29, 5).Value = "L6173"
$ws.Cells.Item(30, 1).Value = "12:48:04"
$ws.Cells.Item(30, 2).Value = "12:54"
$ws.Cells.Item(30, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(30, 4).Value = 6
$ws.Cells.Item(30, 5).Value = "L6203"
$ws.Cells.Item(31, 1).Value = "13:02:04"
$ws.Cells.Item(31, 2).Value = "13:31"
$ws.Cells.Item(31, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(31, 4).Value = 29
$ws.Cells.Item(31, 5).Value = "L6173"
$ws.Cells.Item(32, 1).Value = "13:49:12"
$ws.Cells.Item(32, 2).Value = "14:09"
$ws.Cells.Item(32, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(32, 4).Value = 20
$ws.Cells.Item(32, 5).Value = "L6173"
$ws.Cells.Item(33, 1).Value = "14:47:31"
$ws.Cells.Item(33, 2).Value = "14:47"
$ws.Cells.Item(33, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = "L6203"
$ws.Cells.Item(34, 1).Value = "14:17:09"
$ws.Cells.Item(34, 2).Value = "14:53"
$ws.Cells.Item(34, 3).Value = "215D_LA PLATA"
$ws.Cells.Item(34, 4).Value = 36
$ws.Cells.Item(34, 5).Value = "L6203"
$ws.Cells.Item(35, 1).Value = "14:17:09"
$ws.Cells.Item(35, 2).Value = "15:34"
$ws.Cells.Item(35, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(35, 4).Value = 77
$ws.Cells.Item(35, 5).Value = "L6173"
$ws.Cells.Item(36, 1).Value = "14:47:31"
$ws.Cells.Item(36, 2).Value = "15:38"
$ws.Cells.Item(36, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(36, 4).Value = 51
$ws.Cells.Item(36, 5).Value = "L6173"
$ws.Cells.Item(37, 1).Value = "14:47:31"
$ws.Cells.Item(37, 2).Value = "16:14"
$ws.Cells.Item(37, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(37, 4).Value = 87
$ws.Cells.Item(37, 5).Value = "L6203"

